$d = $word.ActiveDocument

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-ParagraphXmlByText {
    param(
        [string]$needle,
        [string]$innerP
    )
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            $fullXml = $xmlHeader + $innerP + $xmlFooter
            $p.Range.InsertXML($fullXml)
            return $true
        }
    }
    throw "Could not find paragraph containing: $needle"
}

# --- Paragraph: "We strive to allow the community..." -> "The self funding mechanism..." ---
$innerP1 = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/><w:b/><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">The self funding mechanism that pays for development and encourages the community to decide how budget funds are spent.</w:t></w:r></w:p>'
[void](Set-ParagraphXmlByText "We strive to" $innerP1)

# --- Paragraph: "SmartCash mining prevents..." -> "SmartCash now has SmartMining..." ---
$innerP2 = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t xml:space="preserve">SmartCash now has SmartMining that prevents mining attacks. Mining can be done by anyone with a computer with one or more graphics cards. ASICs have yet to be created for the Keccak mining algorithm and it’s probably safe to assume no ASICs will be created for quite some time.</w:t></w:r></w:p>'
[void](Set-ParagraphXmlByText "SmartCash mining prevents" $innerP2)

Write-Output "done"
